# Week 6 solution fix: correct the misspelled author name "J.R.R Tolkien"
# (missing period after the third "R") to the correct "J.R.R. Tolkien"
# everywhere it appears in the Books table (column B - Author).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("J.R.R Tolkien", "J.R.R. Tolkien")

# Cosmetic: leave the view on the last rows of the table, matching the
# author's final on-screen selection when the fix was saved.
$ws.Activate()
$ws.Range("A77").Select()
